$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 18497.1
$ws.Range("J43").Value = 19441.223
$ws.Range("L43").Value = 19441.223
$ws.Range("N43").Value = -19579.223

$ws.Range("H45").Value = 3000
$ws.Range("J45").Value = 3000
$ws.Range("L45").Value = 9000
$ws.Range("N45").Value = -9384

$ws.Range("H51").Value = 16666.334
$ws.Range("J51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -30968

$ws.Range("H135").Value = 4780.5884
$ws.Range("J135").Value = 7639
$ws.Range("L135").Value = 68751
$ws.Range("N135").Value = -73821

$ws.Range("H137").Value = 8148.7095
$ws.Range("I137").Value = 3081.6287
$ws.Range("J137").Value = 14717.148
$ws.Range("K137").Value = 9244.8861
$ws.Range("L137").Value = 44151.444
$ws.Range("M137").Value = -6694.8861
$ws.Range("N137").Value = -49251.444

$ws.Range("H138").Value = 3582.8164
$ws.Range("I138").Value = 5260.75
$ws.Range("J138").Value = 3038.6216
$ws.Range("K138").Value = 15782.25
$ws.Range("L138").Value = 9115.864799999999
$ws.Range("M138").Value = -10642.25
$ws.Range("N138").Value = -19395.8648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5276087
$ws.Range("I32").Value = 7580257.5
$ws.Range("J32").Value = 32112.55
$ws.Range("K32").Value = 7580257.5
$ws.Range("L32").Value = 32112.55
$ws.Range("M32").Value = -7579970.5
$ws.Range("N32").Value = -32686.55

$ws.Range("H45").Value = 4263.6787
$ws.Range("I45").Value = 3935.36
$ws.Range("K45").Value = 3935.36
$ws.Range("M45").Value = -3558.36

$ws.Range("H61").Value = 18960.773
$ws.Range("I61").Value = 13845.667
$ws.Range("J61").Value = 22502
$ws.Range("K61").Value = 13845.667
$ws.Range("L61").Value = 22502
$ws.Range("M61").Value = -13633.667
$ws.Range("N61").Value = -22926

$ws.Range("H110").Value = 9112.429
$ws.Range("I110").Value = 6143.1816
$ws.Range("J110").Value = 19999.666
$ws.Range("K110").Value = 6143.1816
$ws.Range("L110").Value = 19999.666
$ws.Range("M110").Value = -4098.1816
$ws.Range("N110").Value = -24089.666

$ws.Range("H132").Value = 2390385.2
$ws.Range("I132").Value = 4005.75
$ws.Range("K132").Value = 12017.25
$ws.Range("M132").Value = -9487.25

$ws.Range("H136").Value = 18960.773
$ws.Range("I136").Value = 13845.667
$ws.Range("J136").Value = 22502
$ws.Range("K136").Value = 41537.001
$ws.Range("L136").Value = 67506
$ws.Range("M136").Value = -38987.001
$ws.Range("N136").Value = -72606

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 24536.342
$ws.Range("I20").Value = 7156.25
$ws.Range("J20").Value = 37176.41
$ws.Range("K20").Value = 7156.25
$ws.Range("L20").Value = 37176.41
$ws.Range("M20").Value = -6909.25
$ws.Range("N20").Value = -37670.41

$ws.Range("H94").Value = 2424.1714
$ws.Range("I94").Value = 1083.091
$ws.Range("J94").Value = 4693.6924
$ws.Range("K94").Value = 1083.091
$ws.Range("L94").Value = 4693.6924
$ws.Range("M94").Value = -632.0909999999999
$ws.Range("N94").Value = -5595.6924

$ws.Range("H134").Value = 9015.388999999999
$ws.Range("I134").Value = 2237.3
$ws.Range("K134").Value = 6711.900000000001
$ws.Range("M134").Value = -4176.900000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24670.043
$ws.Range("I31").Value = 10246.333
$ws.Range("J31").Value = 40405
$ws.Range("K31").Value = 10246.333
$ws.Range("L31").Value = 40405
$ws.Range("M31").Value = -9951.333000000001
$ws.Range("N31").Value = -40995

$ws.Range("H34").Value = 24670.043
$ws.Range("I34").Value = 10246.333
$ws.Range("J34").Value = 40405
$ws.Range("K34").Value = 10246.333
$ws.Range("L34").Value = 40405
$ws.Range("M34").Value = -10044.333
$ws.Range("N34").Value = -40809

$ws.Range("H132").Value = 9523.913
$ws.Range("J132").Value = 12843.429
$ws.Range("L132").Value = 38530.287
$ws.Range("N132").Value = -43590.287

$ws.Range("H135").Value = 279000
$ws.Range("J135").Value = 279000
$ws.Range("L135").Value = 279000
$ws.Range("N135").Value = -289140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1094.8
$ws.Range("I26").Value = 1094.8
$ws.Range("K26").Value = 3284.4
$ws.Range("M26").Value = -2996.4

$ws.Range("H32").Value = 14068700
$ws.Range("J32").Value = 11112211
$ws.Range("L32").Value = 33336633
$ws.Range("N32").Value = -33337199

$ws.Range("H128").Value = 192473
$ws.Range("I128").Value = 192473
$ws.Range("K128").Value = 577419
$ws.Range("M128").Value = -572439

$ws.Range("H134").Value = 4194.7886
$ws.Range("I134").Value = 1162.6471
$ws.Range("K134").Value = 3487.9413
$ws.Range("M134").Value = 1582.0587

$ws.Range("H137").Value = 2763.913
$ws.Range("I137").Value = 2804.7778
$ws.Range("K137").Value = 8414.3334
$ws.Range("M137").Value = -3314.3334

$ws.Range("H139").Value = 6254.4165
$ws.Range("I139").Value = 5485.5713
$ws.Range("K139").Value = 16456.7139
$ws.Range("M139").Value = -11316.7139

$ws.Range("H140").Value = 1511.5
$ws.Range("I140").Value = 646.1053000000001
$ws.Range("K140").Value = 1938.3159
$ws.Range("M140").Value = 3241.6841

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 9952
$ws.Range("I18").Value = 5000
$ws.Range("J18").Value = 19856
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 19856
$ws.Range("M18").Value = -4707
$ws.Range("N18").Value = -20442

$ws.Range("H29").Value = 3869.1667
$ws.Range("I29").Value = 1751.75
$ws.Range("J29").Value = 8104
$ws.Range("K29").Value = 1751.75
$ws.Range("L29").Value = 8104
$ws.Range("M29").Value = -1461.75
$ws.Range("N29").Value = -8684

$ws.Range("H102").Value = 2958.0356
$ws.Range("I102").Value = 1121.9166
$ws.Range("K102").Value = 1121.9166
$ws.Range("M102").Value = 500.0834

$ws.Range("H113").Value = 102506.22
$ws.Range("I113").Value = 114882
$ws.Range("K113").Value = 114882
$ws.Range("M113").Value = -112712

$ws.Range("H132").Value = 11704.65
$ws.Range("I132").Value = 6212.6924
$ws.Range("J132").Value = 21904
$ws.Range("K132").Value = 18638.0772
$ws.Range("L132").Value = 65712
$ws.Range("M132").Value = -16108.0772
$ws.Range("N132").Value = -70772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1665.5385
$ws.Range("I31").Value = 1670.2
$ws.Range("J31").Value = 1650
$ws.Range("K31").Value = 1670.2
$ws.Range("L31").Value = 1650
$ws.Range("M31").Value = -1422.2
$ws.Range("N31").Value = -2146

$ws.Range("H57").Value = 583333300
$ws.Range("I57").Value = 375000000
$ws.Range("J57").Value = 1000000000
$ws.Range("K57").Value = 375000000
$ws.Range("L57").Value = 1000000000
$ws.Range("M57").Value = -374999434
$ws.Range("N57").Value = -1000001132

$ws.Range("H61").Value = 4017.0557
$ws.Range("I61").Value = 1611
$ws.Range("J61").Value = 5548.1816
$ws.Range("K61").Value = 1611
$ws.Range("L61").Value = 5548.1816
$ws.Range("M61").Value = -1409
$ws.Range("N61").Value = -5952.1816

$ws.Range("H82").Value = 8680.308000000001
$ws.Range("I82").Value = 3649.5557
$ws.Range("J82").Value = 19999.5
$ws.Range("K82").Value = 3649.5557
$ws.Range("L82").Value = 19999.5
$ws.Range("M82").Value = -3288.5557
$ws.Range("N82").Value = -20721.5

$ws.Range("H85").Value = 8680.308000000001
$ws.Range("I85").Value = 3649.5557
$ws.Range("J85").Value = 19999.5
$ws.Range("K85").Value = 3649.5557
$ws.Range("L85").Value = 19999.5
$ws.Range("M85").Value = -2401.5557
$ws.Range("N85").Value = -22495.5

$ws.Range("H113").Value = 4017.0557
$ws.Range("I113").Value = 1611
$ws.Range("J113").Value = 5548.1816
$ws.Range("K113").Value = 1611
$ws.Range("L113").Value = 5548.1816
$ws.Range("M113").Value = 559
$ws.Range("N113").Value = -9888.1816

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H132").Value = 3661455.2
$ws.Range("I132").Value = 7999.5
$ws.Range("K132").Value = 23998.5
$ws.Range("M132").Value = -21468.5

$ws.Range("H136").Value = 13930.439
$ws.Range("I136").Value = 13146.772
$ws.Range("J136").Value = 14837.842
$ws.Range("K136").Value = 39440.31600000001
$ws.Range("L136").Value = 44513.526
$ws.Range("M136").Value = -36890.31600000001
$ws.Range("N136").Value = -49613.526

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 14795.714
$ws.Range("I48").Value = 14828.333
$ws.Range("K48").Value = 14828.333
$ws.Range("M48").Value = -14259.333

$ws.Range("H136").Value = 22740.477
$ws.Range("I136").Value = 982.8889
$ws.Range("K136").Value = 2948.6667
$ws.Range("M136").Value = -398.6667000000002
